$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.410626908494325
$ws.Range("K2").Value = 0.32098402394466
$ws.Range("L2").Value = 0.389807639671357
$ws.Range("N2").Value = 0.400110573948723
